$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 114: correct date + close price (R script re-ran with fresh data) ---
$ws.Range("A114").Value = 45453.2916666667
$ws.Range("C114").Value = 2.6800000667572

# --- Row 115: newly appended observation ---
$ws.Range("A115").Value = 45455.6494675926
$ws.Range("B115").Value = 4000
$ws.Range("C115").Value = 2.59999990463257
$ws.Range("D115").Value = 2.55999994277954
$ws.Range("E115").Value = 2.5699999332428
$ws.Range("F115").Value = 2.58999991416931

# G115 (adj_close) is stored as text in this sheet (matches the rest of
# column G), so build it as a text formula then flatten to a literal value
# via copy / paste-special so no new number-format style gets allocated.
$ws.Range("G115").Formula = "=""2.58999991416931"""
$ws.Range("G115").Copy()
$ws.Range("G115").PasteSpecial(-4163)

$ws.Range("H115").Value = "LS.MI"

# match the date-serial number-format style used by the rest of column A
# (copy format only, so no new style entry is allocated in styles.xml)
$ws.Range("A114").Copy()
$ws.Range("A115").PasteSpecial(-4122)
$excel.CutCopyMode = 0
